# Add a reference to Jester on GitHub in the "Open-Source Work" section.
#
# Before: "I am currently working on Jester - RSpec style testing framework for PowerShell. "
# After:  "I am currently working on Jester(<hyperlink to https://github.com/mbergal/Jester>) - RSpec style testing framework for PowerShell. "

$d = $word.ActiveDocument

# Step 1: turn "Jester - RSpec" into "Jester(<url>) - RSpec" (plain text for now,
# the URL text will be converted into a real hyperlink in step 2).
$d.Content.Find.Execute(
    "Jester - RSpec", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Jester(https://github.com/mbergal/Jester) - RSpec", 2) | Out-Null

# Step 2: locate the URL text we just inserted and convert it into a hyperlink,
# mirroring <w:hyperlink> / <w:rStyle w:val="Hyperlink"/> in the target markup.
$urlRange = $d.Content
$urlRange.Find.Execute(
    "https://github.com/mbergal/Jester", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

$d.Hyperlinks.Add($urlRange, "https://github.com/mbergal/Jester") | Out-Null
